$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'26.969.73"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  -0.55%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'1.825.75"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +0.20%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  -0.35%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'312.60"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +0.15%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('E6').Value = "'  -0.32%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'0.4561"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  -1.42%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = "'  +1.86%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('E9').Value = "'  +0.14%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'0.8707"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  +0.01%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'2.016.73"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  +10.02%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'0.07962"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  +4.29%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'19.71"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  -1.86%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('B14').Value = "'Polkadot"
$ws.Range('B14').Style = 'Normal'
$ws.Range('C14').Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range('C14').Style = 'Normal'
$ws.Range('D14').Value = "'5.314"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  -0.51%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('B15').Value = "'Chainlink"
$ws.Range('B15').Style = 'Normal'
$ws.Range('C15').Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range('C15').Style = 'Normal'
$ws.Range('D15').Value = "'6.520"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +0.65%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'91.38"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  -1.12%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'1.009"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  -0.11%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'0.000008845"
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').Value = "'1.005"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  -0.48%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'14.69"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  +1.49%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'26.733.36"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  -2.48%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'5.102"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  -2.14%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'10.51"
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').Value = "'1.985.39"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  -5.27%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'153.30"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  +1.21%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'1.842"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  -1.58%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'18.31"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  +0.55%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'2.038"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  -1.75%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'5.142"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  +0.90%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'115.12"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  -0.87%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'0.08876"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  -0.37%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D33').Value = "'0.7267"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  -1.19%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'4.405"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  -1.11%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'1.128"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  -0.83%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'1.073"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  +0.17%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'2.442"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  -0.98%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'0.01935"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  +1.13%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D40').Value = "'2.935"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  +0.39%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'7.097"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  -0.73%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('B42').Value = "'Frax"
$ws.Range('B42').Style = 'Normal'
$ws.Range('C42').Value = "'https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range('C42').Style = 'Normal'
$ws.Range('D42').Value = "'0.8915"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  -11.77%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('B43').Value = "'TheSandbox"
$ws.Range('B43').Style = 'Normal'
$ws.Range('C43').Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range('C43').Style = 'Normal'
$ws.Range('D43').Value = "'0.5129"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  -1.34%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('B44').Value = "'Algorand"
$ws.Range('B44').Style = 'Normal'
$ws.Range('C44').Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range('C44').Style = 'Normal'
$ws.Range('D44').Value = "'0.1625"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  -0.10%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('B45').Value = "'Aptos"
$ws.Range('B45').Style = 'Normal'
$ws.Range('C45').Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range('C45').Style = 'Normal'
$ws.Range('D45').Value = "'8.168"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  -1.29%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('B46').Value = "'Decentraland"
$ws.Range('B46').Style = 'Normal'
$ws.Range('C46').Value = "'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range('C46').Style = 'Normal'
$ws.Range('D46').Value = "'0.4820"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -0.32%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('B47').Value = "'PaxDollar"
$ws.Range('B47').Style = 'Normal'
$ws.Range('C47').Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range('C47').Style = 'Normal'
$ws.Range('D47').Value = "'1.007"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  -0.27%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('B48').Value = "'EnergySwap"
$ws.Range('B48').Style = 'Normal'
$ws.Range('C48').Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range('C48').Style = 'Normal'
$ws.Range('D48').Value = "'10.21"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  +0.54%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('B49').Value = "'Quant"
$ws.Range('B49').Style = 'Normal'
$ws.Range('C49').Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range('C49').Style = 'Normal'
$ws.Range('D49').Value = "'102.15"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  -1.19%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('B50').Value = "'NEARProtocol"
$ws.Range('B50').Style = 'Normal'
$ws.Range('C50').Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range('C50').Style = 'Normal'
$ws.Range('D50').Value = "'1.628"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  -0.46%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('B51').Value = "'Cronos"
$ws.Range('B51').Style = 'Normal'
$ws.Range('C51').Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range('C51').Style = 'Normal'
$ws.Range('D51').Value = "'0.06202"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  -0.96%  "
$ws.Range('E51').Style = 'Normal'
